$wb = $excel.ActiveWorkbook

# Rename the "Update_Nightly_RETAIL" / "Update_Nightly_IT" workflow references
# to "Update_Nightly_RE" across the relevant sheets.

$wsWorkflow = $wb.Worksheets.Item("Workflow")
$wsWorkflow.Range("B2").Value = "Update_Nightly_RE"

$wsSteps = $wb.Worksheets.Item("Steps")
$wsSteps.Range("A2").Value = "Update_Nightly_RE"
$wsSteps.Range("A3").Value = "Update_Nightly_RE"

$wsParameters = $wb.Worksheets.Item("Parameters")
$wsParameters.Range("A2").Value = "Update_Nightly_RE"
$wsParameters.Range("A3").Value = "Update_Nightly_RE"

# Best-fit the first column on Parameters following the text edit
$wsParameters.Columns.Item(1).AutoFit() | Out-Null

# Update selections / active cells left over from the editing session
$wsSteps.Range("B6").Select() | Out-Null
$wsParameters.Range("D3").Select() | Out-Null

$wsWorkflow.Activate() | Out-Null
$wsWorkflow.Range("B6").Select() | Out-Null

$wb.Save()
